$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Multiply the values in columns C..M, rows 2..54 by 10,000,000
$rng = $ws.Range("C2:M54")
foreach ($cell in $rng.Cells) {
    $v = $cell.Value()
    if ($v -ne $null) {
        $cell.Value = $v * 10000000
    }
}

# Set column C width to best-fit (auto-fit) per the diff (width ~14.89, bestFit)
$ws.Columns.Item(3).AutoFit() | Out-Null

# Update the view: remove frozen/scrolled topLeftCell, set selection to P16
$ws.Range("P16").Select() | Out-Null
